$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.222587333333333
$ws.Range("H2").Value = 6.667762
$ws.Range("I2").Value = 0.1134117015526119
$ws.Range("J2").Value = 0.1134117015526119
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 87.94127933333334
$ws.Range("N2").Value = 263.823838
$ws.Range("O2").Value = 0.4109331243514438
$ws.Range("P2").Value = 0.4109331243514437
$ws.Range("Q2").Value = 195.4571735233951
$ws.Range("R2").Value = 1759.114561710556
$ws.Range("S2").Value = 0.04660462485702831
$ws.Range("T2").Value = 0.04660462485702831

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.222587333333333
$ws.Range("H3").Value = 6.667762
$ws.Range("I3").Value = 0.1134117015526119
$ws.Range("J3").Value = 0.1134117015526119
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 52.441971
$ws.Range("N3").Value = 157.325913
$ws.Range("O3").Value = 0.2450515065683088
$ws.Range("P3").Value = 0.2450515065683087
$ws.Range("Q3").Value = 116.556860479634
$ws.Range("R3").Value = 1049.011744316706
$ws.Range("S3").Value = 0.02779170832794296
$ws.Range("T3").Value = 0.02779170832794295

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.222587333333333
$ws.Range("H4").Value = 6.667762
$ws.Range("I4").Value = 0.1134117015526119
$ws.Range("J4").Value = 0.1134117015526119
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 54.667459
$ws.Range("N4").Value = 164.002377
$ws.Range("O4").Value = 0.255450795093328
$ws.Range("P4").Value = 0.255450795093328
$ws.Range("Q4").Value = 121.5032019189193
$ws.Range("R4").Value = 1093.528817270274
$ws.Range("S4").Value = 0.02897110933450195
$ws.Range("T4").Value = 0.02897110933450193

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.222587333333333
$ws.Range("H5").Value = 6.667762
$ws.Range("I5").Value = 0.1134117015526119
$ws.Range("J5").Value = 0.1134117015526119
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 18.95316166666667
$ws.Range("N5").Value = 56.85948500000001
$ws.Range("O5").Value = 0.08856457398691947
$ws.Range("P5").Value = 0.08856457398691944
$ws.Range("Q5").Value = 42.12505704695223
$ws.Range("R5").Value = 379.12551342257
$ws.Range("S5").Value = 0.01004425903313873
$ws.Range("T5").Value = 0.01004425903313873

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.89424733333333
$ws.Range("H6").Value = 35.682742
$ws.Range("I6").Value = 0.6069263549423107
$ws.Range("J6").Value = 0.6069263549423106
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 87.94127933333334
$ws.Range("N6").Value = 263.823838
$ws.Range("O6").Value = 0.4109331243514438
$ws.Range("P6").Value = 0.4109331243514437
$ws.Range("Q6").Value = 1045.995327200422
$ws.Range("R6").Value = 9413.957944803797
$ws.Range("S6").Value = 0.2494061432876771
$ws.Range("T6").Value = 0.249406143287677

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.89424733333333
$ws.Range("H7").Value = 35.682742
$ws.Range("I7").Value = 0.6069263549423107
$ws.Range("J7").Value = 0.6069263549423106
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 52.441971
$ws.Range("N7").Value = 157.325913
$ws.Range("O7").Value = 0.2450515065683088
$ws.Range("P7").Value = 0.2450515065683087
$ws.Range("Q7").Value = 623.757773721494
$ws.Range("R7").Value = 5613.819963493446
$ws.Range("S7").Value = 0.1487282176546253
$ws.Range("T7").Value = 0.1487282176546253

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.89424733333333
$ws.Range("H8").Value = 35.682742
$ws.Range("I8").Value = 0.6069263549423107
$ws.Range("J8").Value = 0.6069263549423106
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 54.667459
$ws.Range("N8").Value = 164.002377
$ws.Range("O8").Value = 0.255450795093328
$ws.Range("P8").Value = 0.255450795093328
$ws.Range("Q8").Value = 650.2282784308593
$ws.Range("R8").Value = 5852.054505877733
$ws.Range("S8").Value = 0.1550398199331087
$ws.Range("T8").Value = 0.1550398199331086

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.89424733333333
$ws.Range("H9").Value = 35.682742
$ws.Range("I9").Value = 0.6069263549423107
$ws.Range("J9").Value = 0.6069263549423106
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 18.95316166666667
$ws.Range("N9").Value = 56.85948500000001
$ws.Range("O9").Value = 0.08856457398691947
$ws.Range("P9").Value = 0.08856457398691944
$ws.Range("Q9").Value = 225.4335926119856
$ws.Range("R9").Value = 2028.90233350787
$ws.Range("S9").Value = 0.05375217406689962
$ws.Range("T9").Value = 0.0537521740668996

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.762
$ws.Range("H10").Value = 2.286
$ws.Range("I10").Value = 0.0388824840702579
$ws.Range("J10").Value = 0.03888248407025789
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 87.94127933333334
$ws.Range("N10").Value = 263.823838
$ws.Range("O10").Value = 0.4109331243514438
$ws.Range("P10").Value = 0.4109331243514437
$ws.Range("Q10").Value = 67.01125485200001
$ws.Range("R10").Value = 603.101293668
$ws.Range("S10").Value = 0.01597810066153632
$ws.Range("T10").Value = 0.01597810066153632

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.762
$ws.Range("H11").Value = 2.286
$ws.Range("I11").Value = 0.0388824840702579
$ws.Range("J11").Value = 0.03888248407025789
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 52.441971
$ws.Range("N11").Value = 157.325913
$ws.Range("O11").Value = 0.2450515065683088
$ws.Range("P11").Value = 0.2450515065683087
$ws.Range("Q11").Value = 39.960781902
$ws.Range("R11").Value = 359.647037118
$ws.Range("S11").Value = 0.009528211300534965
$ws.Range("T11").Value = 0.009528211300534961

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.762
$ws.Range("H12").Value = 2.286
$ws.Range("I12").Value = 0.0388824840702579
$ws.Range("J12").Value = 0.03888248407025789
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 54.667459
$ws.Range("N12").Value = 164.002377
$ws.Range("O12").Value = 0.255450795093328
$ws.Range("P12").Value = 0.255450795093328
$ws.Range("Q12").Value = 41.656603758
$ws.Range("R12").Value = 374.909433822
$ws.Range("S12").Value = 0.009932561470951042
$ws.Range("T12").Value = 0.009932561470951038

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.762
$ws.Range("H13").Value = 2.286
$ws.Range("I13").Value = 0.0388824840702579
$ws.Range("J13").Value = 0.03888248407025789
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 18.95316166666667
$ws.Range("N13").Value = 56.85948500000001
$ws.Range("O13").Value = 0.08856457398691947
$ws.Range("P13").Value = 0.08856457398691944
$ws.Range("Q13").Value = 14.44230919
$ws.Range("R13").Value = 129.98078271
$ws.Range("S13").Value = 0.003443610637235573
$ws.Range("T13").Value = 0.003443610637235572

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.718678666666666
$ws.Range("H14").Value = 14.156036
$ws.Range("I14").Value = 0.2407794594348195
$ws.Range("J14").Value = 0.2407794594348195
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 87.94127933333334
$ws.Range("N14").Value = 263.823838
$ws.Range("O14").Value = 0.4109331243514438
$ws.Range("P14").Value = 0.4109331243514437
$ws.Range("Q14").Value = 414.9666387095742
$ws.Range("R14").Value = 3734.699748386168
$ws.Range("S14").Value = 0.09894425554520209
$ws.Range("T14").Value = 0.09894425554520207

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.718678666666666
$ws.Range("H15").Value = 14.156036
$ws.Range("I15").Value = 0.2407794594348195
$ws.Range("J15").Value = 0.2407794594348195
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 52.441971
$ws.Range("N15").Value = 157.325913
$ws.Range("O15").Value = 0.2450515065683088
$ws.Range("P15").Value = 0.2450515065683087
$ws.Range("Q15").Value = 247.456809795652
$ws.Range("R15").Value = 2227.111288160868
$ws.Range("S15").Value = 0.05900336928520549
$ws.Range("T15").Value = 0.05900336928520548

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.718678666666666
$ws.Range("H16").Value = 14.156036
$ws.Range("I16").Value = 0.2407794594348195
$ws.Range("J16").Value = 0.2407794594348195
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 54.667459
$ws.Range("N16").Value = 164.002377
$ws.Range("O16").Value = 0.255450795093328
$ws.Range("P16").Value = 0.255450795093328
$ws.Range("Q16").Value = 257.9581725441747
$ws.Range("R16").Value = 2321.623552897572
$ws.Range("S16").Value = 0.06150730435476635
$ws.Range("T16").Value = 0.06150730435476634

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.718678666666666
$ws.Range("H17").Value = 14.156036
$ws.Range("I17").Value = 0.2407794594348195
$ws.Range("J17").Value = 0.2407794594348195
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 18.95316166666667
$ws.Range("N17").Value = 56.85948500000001
$ws.Range("O17").Value = 0.08856457398691947
$ws.Range("P17").Value = 0.08856457398691944
$ws.Range("Q17").Value = 89.43387962238445
$ws.Range("R17").Value = 804.9049166014601
$ws.Range("S17").Value = 0.02132453024964554
$ws.Range("T17").Value = 0.02132453024964554
